$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Junio de 2020 a las 04:26"

# --- Brasil (row 5): update Recuperados / Casos activos ---
$ws.Range("D5").Value = 594104
$ws.Range("E5").Value = 465837

# --- Bolivia / Irlanda swap + updated figures ---
# Row 47 becomes "Bolivia" with freshly reported totals.
$ws.Range("A47").Value = "Bolivia"
$ws.Range("B47").Value = 25493
$ws.Range("C47").Value = 1105
$ws.Range("D47").Value = 5857
$ws.Range("E47").Value = 18816
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 47
$ws.Range("H47").Value = 820

# Row 48 becomes "Irlanda", carrying the figures Irlanda previously had.
$ws.Range("A48").Value = "Irlanda"
$ws.Range("B48").Value = 25383
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 22698
$ws.Range("E48").Value = 968
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 1717

# --- Corea del Sur (row 64): figures refresh ---
$ws.Range("B64").Value = 12484
$ws.Range("C64").Value = 46
$ws.Range("D64").Value = 10908
$ws.Range("E64").Value = 1295
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 1
$ws.Range("H64").Value = 281

# --- Venezuela / Bulgaria swap + updated figures ---
# Row 90 becomes "Venezuela" with freshly reported totals.
$ws.Range("A90").Value = "Venezuela"
$ws.Range("B90").Value = 4048
$ws.Range("C90").Value = 0
$ws.Range("D90").Value = 1327
$ws.Range("E90").Value = 2686
$ws.Range("F90").Value = 0
$ws.Range("G90").Value = 0
$ws.Range("H90").Value = 35

# Row 91 becomes "Bulgaria", carrying the figures Bulgaria previously had.
$ws.Range("A91").Value = "Bulgaria"
$ws.Range("B91").Value = 3984
$ws.Range("C91").Value = 0
$ws.Range("D91").Value = 2171
$ws.Range("E91").Value = 1606
$ws.Range("F91").Value = 0
$ws.Range("G91").Value = 0
$ws.Range("H91").Value = 207

# --- Nueva Zelanda (row 119): figures refresh ---
$ws.Range("B119").Value = 1515
$ws.Range("C119").Value = 2
$ws.Range("D119").Value = 1483
$ws.Range("E119").Value = 10
